# This workbook tracks daily/weekly price records for Coliflor at
# "Feria Lagunitas de Puerto Montt". The update adds one new weekly
# record. The new record is inserted as the new row 290 (pushing every
# existing row from 290 downward by one), matching the source diff where
# the dimension grows from A1:R418 to A1:R419 and every row D..418
# (old) becomes D+1..419 (new).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 290; everything currently at row 290
# (and below) shifts down to row 291 (and below).
$ws.Rows.Item(290).Insert()

# Fill in the new weekly record in row 290.
$ws.Range("A290").Value2 = 4
$ws.Range("B290").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C290").Value2 = "Los Lagos"
$ws.Range("D290").Value2 = 44839
$ws.Range("E290").Value2 = 10
$ws.Range("F290").Value2 = 100112008
$ws.Range("G290").Value2 = "Coliflor"
$ws.Range("H290").Value2 = "Sin especificar"
$ws.Range("I290").Value2 = "Primera"
$ws.Range("J290").Value2 = 100
$ws.Range("K290").Value2 = 1800
$ws.Range("L290").Value2 = 1800
$ws.Range("M290").Value2 = 1800
$ws.Range("N290").Value2 = "`$/unidad"
$ws.Range("O290").Value2 = "Región Metropolitana"
$ws.Range("P290").Value2 = 1800
$ws.Range("Q290").Value2 = 1
$ws.Range("R290").Value2 = "Hortaliza"
